$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column F
$ws.Range("F1").Value = "sd_temp"

# Standard deviation values for each row (2-17)
$sdValues = @(
    4.15341136580214,
    3.34061940946596,
    2.92543810392176,
    2.23659585411325,
    3.38340271378673,
    3.84360641907832,
    2.6684827713548,
    3.33033984548803,
    3.077139100033,
    4.54695424496783,
    2.51979465348512,
    3.86815126501919,
    4.44997724283534,
    3.29938177983242,
    2.97227958642577,
    3.75335425017699
)

$row = 2
foreach ($val in $sdValues) {
    $ws.Cells.Item($row, 6).Value = $val
    $row = $row + 1
}
